$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the existing "Residuals NOT correlated..." cell keeps its wrapped style,
# so the new cell below it (E4) can reuse the same cell style (index) rather than
# the engine minting a brand new one.
$ws.Range("H3").WrapText = $true

# New comparison row: SOL vs BEL BT
$ws.Range("A4").Value = "SOL"
$ws.Range("B4").Value = "BEL BT"
$ws.Range("C4").Value = 1.729352
$ws.Range("D4").Value = -0.006369
$ws.Range("E4").Value = 0.455329
$ws.Range("F4").Value = 0.265
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = -92.00044
$ws.Range("I4").Value = "No"
$ws.Range("J4").Value = "No"
$ws.Range("K4").Value = "No"
$ws.Range("L4").Value = "No"

# E4 shares the wrapped-text style used by H3
$ws.Range("E4").WrapText = $true

# Row 4 renders at the same (slightly taller) height as row 3
$ws.Rows.Item(4).RowHeight = 14.9

# Move the active selection past the newly added row
$ws.Range("L5").Select()
